$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ref")
$ws.Range("A1").Value = "TEST"
